# Apply the "Add data for 2021-10-22" update to the carjacking arrests workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date.
$ws.Name = "Through 2021-10-14"

# Row 6 (April) - 2021 column group (T/U/V)
$ws.Range("T6").Value = 11
$ws.Range("U6").Value = 89
$ws.Range("V6").Value = 0.11

# Row 12 (October) - update label and 2015-2020 groups
$ws.Range("A12").Value = "October (through 10-14)"

$ws.Range("C12").Value = 12
$ws.Range("D12").Value = 0.0769

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 20
$ws.Range("G12").Value = 0.0909

$ws.Range("I12").Value = 18
$ws.Range("J12").Value = 0.2174

$ws.Range("L12").Value = 34
$ws.Range("M12").Value = 0.0556

$ws.Range("O12").Value = 18
$ws.Range("P12").Value = 0.0526

$ws.Range("R12").Value = 65
$ws.Range("U12").Value = 86

# Row 13 (Total)
$ws.Range("C13").Value = 208
$ws.Range("D13").Value = 0.1297

$ws.Range("E13").Value = 48
$ws.Range("F13").Value = 403
$ws.Range("G13").Value = 0.1064

$ws.Range("I13").Value = 595
$ws.Range("J13").Value = 0.0846

$ws.Range("L13").Value = 521
$ws.Range("M13").Value = 0.1079

$ws.Range("O13").Value = 397
$ws.Range("P13").Value = 0.0998

$ws.Range("R13").Value = 913
$ws.Range("S13").Value = 0.0549

$ws.Range("T13").Value = 82
$ws.Range("U13").Value = 1254
$ws.Range("V13").Value = 0.0614
